$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 200.6923
$ws.Range("I6").Value = 209.08333
$ws.Range("K6").Value = 627.24999
$ws.Range("M6").Value = -515.24999
$ws.Range("H8").Value = 60
$ws.Range("I8").Value = 60
$ws.Range("K8").Value = 180
$ws.Range("M8").Value = -41
$ws.Range("H51").Value = 4999.4116
$ws.Range("J51").Value = 4999.4116
$ws.Range("L51").Value = 4999.4116
$ws.Range("N51").Value = -5967.4116
$ws.Range("H112").Value = 5884537.5
$ws.Range("J112").Value = 5954576
$ws.Range("L112").Value = 17863728
$ws.Range("N112").Value = -17865944
$ws.Range("H138").Value = 8337984
$ws.Range("I138").Value = 1056.6111
$ws.Range("K138").Value = 3169.8333
$ws.Range("M138").Value = 1970.1667

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 650.25
$ws.Range("I5").Value = 650.25
$ws.Range("K5").Value = 650.25
$ws.Range("M5").Value = -538.25
$ws.Range("H35").Value = 2491.6667
$ws.Range("I35").Value = 1987.5
$ws.Range("K35").Value = 1987.5
$ws.Range("M35").Value = -1581.5
$ws.Range("H61").Value = 38464936
$ws.Range("I61").Value = 55556988
$ws.Range("K61").Value = 55556988
$ws.Range("M61").Value = -55556776
$ws.Range("H63").Value = 4145.108
$ws.Range("I63").Value = 2800.7144
$ws.Range("J63").Value = 8327.666999999999
$ws.Range("K63").Value = 2800.7144
$ws.Range("L63").Value = 8327.666999999999
$ws.Range("M63").Value = -2114.7144
$ws.Range("N63").Value = -9699.666999999999
$ws.Range("H66").Value = 4145.108
$ws.Range("I66").Value = 2800.7144
$ws.Range("J66").Value = 8327.666999999999
$ws.Range("K66").Value = 14003.572
$ws.Range("L66").Value = 41638.335
$ws.Range("M66").Value = -10571.572
$ws.Range("N66").Value = -48502.335
$ws.Range("H74").Value = 66743144
$ws.Range("I74").Value = 77010450
$ws.Range("K74").Value = 77010450
$ws.Range("M74").Value = -77009576
$ws.Range("H77").Value = 66743144
$ws.Range("I77").Value = 77010450
$ws.Range("K77").Value = 385052250
$ws.Range("M77").Value = -385047882
$ws.Range("H132").Value = 26317870
$ws.Range("I132").Value = 1958.8286
$ws.Range("K132").Value = 5876.4858
$ws.Range("M132").Value = -3346.4858
$ws.Range("H136").Value = 38464936
$ws.Range("I136").Value = 55556988
$ws.Range("K136").Value = 166670964
$ws.Range("M136").Value = -166668414

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 650.25
$ws.Range("I4").Value = 650.25
$ws.Range("K4").Value = 650.25
$ws.Range("M4").Value = -535.25
$ws.Range("H86").Value = 12325.807
$ws.Range("I86").Value = 6435.4585
$ws.Range("K86").Value = 6435.4585
$ws.Range("M86").Value = -5312.4585
$ws.Range("H89").Value = 12325.807
$ws.Range("I89").Value = 6435.4585
$ws.Range("K89").Value = 32177.2925
$ws.Range("M89").Value = -26561.2925
$ws.Range("H99").Value = 3400.5881
$ws.Range("I99").Value = 1659.2727
$ws.Range("K99").Value = 1659.2727
$ws.Range("M99").Value = -161.2727

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 401.57144
$ws.Range("I2").Value = 342.4
$ws.Range("J2").Value = 549.5
$ws.Range("K2").Value = 342.4
$ws.Range("L2").Value = 549.5
$ws.Range("M2").Value = -229.4
$ws.Range("N2").Value = -775.5
$ws.Range("H56").Value = 51398.8
$ws.Range("I56").Value = 19664.666
$ws.Range("J56").Value = 99000
$ws.Range("K56").Value = 19664.666
$ws.Range("L56").Value = 99000
$ws.Range("M56").Value = -18819.666
$ws.Range("N56").Value = -100690

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 25000
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 25000
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 25000
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -25554
$ws.Range("H97").Value = 1585.8125
$ws.Range("I97").Value = 484.625
$ws.Range("K97").Value = 484.625
$ws.Range("M97").Value = 11.375
$ws.Range("H113").Value = 3148.5217
$ws.Range("J113").Value = 4690.2
$ws.Range("L113").Value = 4690.2
$ws.Range("N113").Value = -9030.200000000001
$ws.Range("H132").Value = 2853.9
$ws.Range("I132").Value = 2844.72
$ws.Range("K132").Value = 8534.16
$ws.Range("M132").Value = -6004.16

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1128.8846
$ws.Range("I46").Value = 576.7907
$ws.Range("J46").Value = 3766.6667
$ws.Range("K46").Value = 576.7907
$ws.Range("L46").Value = 3766.6667
$ws.Range("M46").Value = -388.7907
$ws.Range("N46").Value = -4142.6667
$ws.Range("H61").Value = 4697.923
$ws.Range("I61").Value = 3712.5715
$ws.Range("K61").Value = 3712.5715
$ws.Range("M61").Value = -3510.5715
$ws.Range("H82").Value = 2840.1875
$ws.Range("I82").Value = 2083.4443
$ws.Range("J82").Value = 3813.1428
$ws.Range("K82").Value = 2083.4443
$ws.Range("L82").Value = 3813.1428
$ws.Range("M82").Value = -1722.4443
$ws.Range("N82").Value = -4535.1428
$ws.Range("H85").Value = 2840.1875
$ws.Range("I85").Value = 2083.4443
$ws.Range("J85").Value = 3813.1428
$ws.Range("K85").Value = 2083.4443
$ws.Range("L85").Value = 3813.1428
$ws.Range("M85").Value = -835.4443000000001
$ws.Range("N85").Value = -6309.1428
$ws.Range("H100").Value = 4162.625
$ws.Range("I100").Value = 2799.6667
$ws.Range("J100").Value = 4980.4
$ws.Range("K100").Value = 2799.6667
$ws.Range("L100").Value = 4980.4
$ws.Range("M100").Value = -2258.6667
$ws.Range("N100").Value = -6062.4
$ws.Range("H113").Value = 4697.923
$ws.Range("I113").Value = 3712.5715
$ws.Range("K113").Value = 3712.5715
$ws.Range("M113").Value = -1542.5715
$ws.Range("H122").Value = 4675.5
$ws.Range("I122").Value = 3888.25
$ws.Range("K122").Value = 11664.75
$ws.Range("M122").Value = -9214.75
$ws.Range("H136").Value = 2264.25
$ws.Range("I136").Value = 1977.9259
$ws.Range("K136").Value = 5933.7777
$ws.Range("M136").Value = -3383.7777

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 75000
$ws.Range("J27").Value = 75000
$ws.Range("L27").Value = 75000
$ws.Range("N27").Value = -75138
$ws.Range("H62").Value = 7570.7856
$ws.Range("I62").Value = 7200.5
$ws.Range("K62").Value = 7200.5
$ws.Range("M62").Value = -6576.5
$ws.Range("H65").Value = 7570.7856
$ws.Range("I65").Value = 7200.5
$ws.Range("K65").Value = 36002.5
$ws.Range("M65").Value = -32882.5
$ws.Range("H96").Value = 7105.2
$ws.Range("I96").Value = 5191.1665
$ws.Range("J96").Value = 9976.25
$ws.Range("K96").Value = 5191.1665
$ws.Range("L96").Value = 9976.25
$ws.Range("M96").Value = -3818.1665
$ws.Range("N96").Value = -12722.25
$ws.Range("H107").Value = 476.26666
$ws.Range("I107").Value = 443.22223
$ws.Range("J107").Value = 525.8333
$ws.Range("K107").Value = 1329.66669
$ws.Range("L107").Value = 1577.4999
$ws.Range("M107").Value = 590.33331
$ws.Range("N107").Value = -5417.4999
$ws.Range("H115").Value = 83999.5
$ws.Range("J115").Value = 83999.5
$ws.Range("L115").Value = 83999.5
$ws.Range("N115").Value = -87133.5
$ws.Range("H122").Value = 83418540
$ws.Range("I122").Value = 91001630
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 273004890
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -273002440
$ws.Range("N122").Value = -18400
$ws.Range("H132").Value = 5292.475
$ws.Range("I132").Value = 5540.875
$ws.Range("K132").Value = 16622.625
$ws.Range("M132").Value = -14092.625
